$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) for rows 2-49 ---
# NOTE: several Price values are plain 'digits.digits' strings (e.g. '551.90').
# Assigning those to .Value directly would let Excel auto-detect them as
# numbers and silently drop the trailing zero (551.90 -> 551.9), so for those
# cells we flip NumberFormat to Text ('@') for the write, then back to
# General so no explicit formatting is left behind.
$ws.Range("D2").Value = "59.488.28"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.572.79"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.90"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.89"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "2.578.99"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.162"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +7.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "3.028.18"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "59.451.16"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.08"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +4.49%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "2.575.01"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.13"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.478"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +6.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.69"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "0.0₃0770"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.68"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.07"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.898"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.52"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.851"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.83"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -3.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.34"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +7.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0969"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.591"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.68"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0233"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.44%  "

# --- Rows 50 & 51 swap ranking: Maker moves up to 50th, InjectiveProtocol to 51st ---
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.970.40"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.69"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.97%  "
